# Logged Week 15 and simulated Week 16
# Update the "R" (Road) row (row 3) target-depth counting stats on both the
# OFF and DEF sheets of the Cardinals 2021 Target Depth Data workbook.

$wb = $excel.ActiveWorkbook

# --- OFF sheet: Short Att, Short Comp, Deep Att, Deep Comp, Short Int ---
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 206
$wsOff.Range("C3").Value = 153
$wsOff.Range("D3").Value = 50
$wsOff.Range("E3").Value = 32
$wsOff.Range("F3").Value = 1

# --- DEF sheet: Short Att, Short Comp, Deep Att, Deep Comp ---
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 164
$wsDef.Range("C3").Value = 112
$wsDef.Range("D3").Value = 30
$wsDef.Range("E3").Value = 14
